$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 254; this shifts existing rows 254-333 down to 255-334
$ws.Rows.Item(254).Insert()

# Populate the newly inserted row 254 with the new data record
$ws.Cells.Item(254, 1).Value = 8
$ws.Cells.Item(254, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(254, 3).Value = "Coquimbo"
$ws.Cells.Item(254, 4).Value = 44876
$ws.Cells.Item(254, 4).NumberFormat = $ws.Cells.Item(255, 4).NumberFormat
$ws.Cells.Item(254, 5).Value = 4
$ws.Cells.Item(254, 6).Value = 100112012
$ws.Cells.Item(254, 7).Value = "Espinaca"
$ws.Cells.Item(254, 8).Value = "Sin especificar"
$ws.Cells.Item(254, 9).Value = "Primera"
$ws.Cells.Item(254, 10).Value = 1600
$ws.Cells.Item(254, 11).Value = 500
$ws.Cells.Item(254, 12).Value = 600
$ws.Cells.Item(254, 13).Value = 550
$ws.Cells.Item(254, 14).Value = "$/atado 300 a 500 gramos"
$ws.Cells.Item(254, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(254, 16).Value = 1100
$ws.Cells.Item(254, 17).Value = 0.5
$ws.Cells.Item(254, 18).Value = "Hortaliza"
